function Replace-Text($doc, $old, $new) {
    $r = $doc.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
        return
    }
    $r.Text = $new
}

$d = $word.ActiveDocument

Replace-Text $d "Appendix 7: SWIFT Focus Group Information Sheet and Consent Form: " "Bylae 7: SWIFT Fokusgroep Inligtingsblad en Toestemmingsvorm: "
Replace-Text $d "Consent Form: Participants" "Toestemmingsvorm: Deelnemers"
Replace-Text $d "Participants" "Deelnemers"
Replace-Text $d "What is a focus group and what will this one look like?" "Wat is 'n fokusgroep en hoe sal hierdie een lyk?"
Replace-Text $d "Why have I been invited to the interview?" "Waarom is ek uitgenooi vir die onderhoud?"
Replace-Text $d "Do I have to agree to be interviewed?" "Moet ek instem om ondervra te word?"
Replace-Text $d "What happens if I agree to be interviewed?" "Wat gebeur as ek instem om ondervra te word?"
Replace-Text $d "Do I get anything for being interviewed? " "Kry ek enige iets vir deelname aan die onderhoud? "
Replace-Text $d "As a thank you for taking part in the discussion, we'll give you a R120 Shoprite voucher afterwards. " "As 'n bedanking vir jou deelname aan die bespreking, sal ons vir jou 'n R120 Shoprite-koepon gee nadat die onderhoud voltooi is. "
Replace-Text $d "What happens to my information if I agree to be interviewed?" "Wat gebeur met my inligting as ek instem om ondervra te word?"
Replace-Text $d "What happens to the research results?" "Wat gebeur met die navorsingsresultate?"
Replace-Text $d "Who are some of the study team members?" "Wie is sommige van die studie-spanlede?"
Replace-Text $d "The principal investigators of this study are Prof Cathy Ward and Cindee Bruyns and the Co-investigator is Carly Katzef all from the University of Cape Town." "Die hoofondersoeker van hierdie studie is Prof. Cathy Ward en Cindee Bruyns, en die Mede-navorser is Carly Katzef, almal van die Universiteit van Kaapstad."
Replace-Text $d "Who pays for the study?" "Wie betaal vir die studie?"
Replace-Text $d "This study is part of the Global Parenting Initiative, funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. " "Hierdie studie is deel van die Global Parenting Initiative, gefinansier deur die LEGO Foundation, Oak Foundation, die World Childhood Foundation, The Human Safety Net, en die UK Research and Innovation Global Challenges Research Fund. "
Replace-Text $d "Data protection" "Databeskerming"
Replace-Text $d "Who has approved this study?" "Wie het hierdie studie goedgekeur?"
Replace-Text $d "Who do I contact if I have questions or concerns?" "Wie kan ek kontak as ek vrae of bekommernisse het?"
Replace-Text $d "If you have any questions or concerns about your rights as a study participant, you can contact the study team at " "As jy enige vrae of bekommernisse het oor jou regte as 'n studie-deelnemer, kan jy die studiespan kontak by "
Replace-Text $d " or on WhatsApp at +27 XX XXX XXXX (messages only)." " of via WhatsApp by +27 XX XXX XXXX (net boodskappe)."
Replace-Text $d "If you have more questions or concerns about your rights, you can contact one of the ethics committees listed: " "As jy meer vrae of bekommernisse het oor jou regte, kan jy een van die etiekkomitees hieronder kontak: "
Replace-Text $d "Name of interviewer" "Naam van onderhoudvoerder"
Replace-Text $d "Name" "Naam"
Replace-Text $d "Telephone" "Telefoon"
Replace-Text $d "Email" "E-pos"
Replace-Text $d "University of Cape Town Centre for Social Science Research " "Universiteit van Kaapstad Sentrum vir Sosiale Wetenskap Navorsing "
Replace-Text $d "Human Research Ethics Committee" "Etiekkomitee vir Menslike Navorsing"
Replace-Text $d "[to be signed at the focus group discussion]" "[moet geteken word by die fokusgroepbespreking]"
Replace-Text $d "I ______________________ (participant name) have read the information above and agree to the following: " "Ek ______________________ (deelnemer naam) het die bogenoemde inligting gelees en stem saam met die volgende: "
Replace-Text $d "Someone from the research team has gone over all the information above and I know what I need to do." "Iemand van die navorsingspan het al die inligting hierbo deurgegaan en ek weet wat ek moet doen."
Replace-Text $d "I know who can see my information after the focus group, how it will be kept safe, and what happens to it after the study." "Ek weet wie my inligting sal kan sien na die fokusgroep, hoe dit veilig gehou sal word, en wat met dit sal gebeur na die studie."
Replace-Text $d "I know I can request access to my data, correct any mistakes, ask to delete it, or for it to be transferred somewhere else." "Ek weet ek kan toegang tot my data versoek, enige foute regstel, vra dat dit verwyder word, of vir dit om na 'n ander plek oorgedra te word."
Replace-Text $d "I know that I won’t be named in any papers or reports from this study." "Ek weet dat ek nie in enige artikels of verslae van hierdie studie genoem sal word nie."
Replace-Text $d "I know who to tell if I have a problem with the study." "Ek weet wie ek moet kontak as ek 'n probleem met die studie het."
Replace-Text $d "I can be contacted again if more information is needed from me." "Ek kan weer gekontak word as meer inligting van my nodig is."
Replace-Text $d "I understand the team will keep my contact information safe so they can tell me about the results of the study." "Ek verstaan dat die span my kontakbesonderhede veilig sal hou sodat hulle my kan inlig oor die resultate van die studie."
Replace-Text $d "If you agree with all the statements above and want to be interviewed for the study, please add your name, signature and the date below." "As jy saamstem met al die bogenoemde stellings en ondervrae wil word vir die studie, voeg asseblief jou naam, handtekening en die datum hieronder by."
Replace-Text $d "Date of Interview" "Datum van Onderhoud"
Replace-Text $d "Signature of Interviewer" "Handtekening van Onderhoudvoerder"
